$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 307
$ws.Range("F4").Value = 415
$ws.Range("F5").Value = 8526
$ws.Range("F7").Value = 10652
$ws.Range("F9").Value = 11
$ws.Range("F10").Value = 18
$ws.Range("F13").Value = 115
$ws.Range("F18").Value = 76
$ws.Range("F22").Value = 1811
$ws.Range("F23").Value = 72
$ws.Range("F24").Value = 543
$ws.Range("F25").Value = 342
$ws.Range("F26").Value = 286
$ws.Range("F27").Value = 61
$ws.Range("F28").Value = 583
$ws.Range("F30").Value = 1176
$ws.Range("F34").Value = 439
$ws.Range("F35").Value = 341
$ws.Range("F36").Value = 284
$ws.Range("F37").Value = 20
$ws.Range("F41").Value = 93
$ws.Range("F42").Value = 530
$ws.Range("F43").Value = 641
$ws.Range("F45").Value = 95
$ws.Range("F46").Value = 93

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 42

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 211
$ws.Range("F3").Value = 2800
$ws.Range("F4").Value = 341

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 307
$ws.Range("F3").Value = 211
$ws.Range("F6").Value = 341
$ws.Range("F9").Value = 415
$ws.Range("F10").Value = 8526
$ws.Range("F12").Value = 10652
$ws.Range("F15").Value = 115
$ws.Range("F17").Value = 76
$ws.Range("F19").Value = 1811
$ws.Range("F20").Value = 72
$ws.Range("F21").Value = 543
$ws.Range("F22").Value = 286
$ws.Range("F23").Value = 61
$ws.Range("F25").Value = 583
$ws.Range("F26").Value = 42
$ws.Range("F28").Value = 1176
$ws.Range("F35").Value = 439
$ws.Range("F37").Value = 341
$ws.Range("F42").Value = 93
$ws.Range("F43").Value = 533
$ws.Range("F47").Value = 641
$ws.Range("F48").Value = 95
$ws.Range("F49").Value = 93
